# Daily attendance processing - 2025-11-17 08:30:06
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the "Recorded By" contributor lists for the two already-recorded
# sessions (ANATOMY session 1 & 2) to reflect the latest processing order.
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg"

# ANATOMY session 3 (row 4) has passed its scheduled time without being
# recorded, so its status moves from "Pending" to "Not Recorded".
$ws.Range("I4").Value = "Not Recorded"

# Update the summary statistics to reflect the status change above:
# one more "Missing" (Not Recorded) session, one fewer "Pending" session.
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 19
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 19
